$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14485.46121151955
$ws.Range("C2").Value = 33734.78440102194
$ws.Range("D2").Value = 69568.33959421814
$ws.Range("E2").Value = 106803.4685720019

$ws.Range("B3").Value = 145003.3954756857
$ws.Range("C3").Value = 298145.8313006964
$ws.Range("D3").Value = 386440.3209392101
$ws.Range("E3").Value = 437253.6901728552

$ws.Range("B4").Value = 15942.47343396684
$ws.Range("C4").Value = 33362.88336293728
$ws.Range("D4").Value = 56608.03133420134
$ws.Range("E4").Value = 74459.8499575834

$ws.Range("B6").Value = 60096.74729507345
$ws.Range("C6").Value = 99153.67958249315
$ws.Range("D6").Value = 106104.9411789463
$ws.Range("E6").Value = 94580.64645259328

$ws.Range("B7").Value = 6747.90035708273
$ws.Range("C7").Value = 16819.5337263134
$ws.Range("D7").Value = 21084.82219611427
$ws.Range("E7").Value = 24848.54708038423

$ws.Range("B9").Value = 504892.9070988327
$ws.Range("C9").Value = 1000724.314930425
$ws.Range("D9").Value = 1540441.836425559
$ws.Range("E9").Value = 2052747.6651996

$ws.Range("B12").Value = 597460.1008527816
$ws.Range("C12").Value = 891326.4376878701
$ws.Range("D12").Value = 881601.2396952431
$ws.Range("E12").Value = 694918.7786644093
